# Change to use rmtr account number for validation
# On the "Test Data (Validate)" sheet, the transfer_unique_no column (G)
# previously held unique dummy values (N01234567890-N01234567895) that were
# used to validate repeated calls. Switch that column to a generic
# "System Generated" placeholder and instead give each row its own unique
# rmtr_account value (column J) so that uniqueness/validation is now keyed
# off the remitter account number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Data (Validate)")

# Column G (transfer_unique_no) rows 8-13: replace per-row unique values
# with a common "System Generated" marker.
$ws.Range("G8:G13").Value = "System Generated"
$ws.Columns.Item(7).AutoFit() | Out-Null

# Column J (rmtr_account) rows 8-13: give each row a unique account number.
$ws.Range("J8").Value = 123456780
$ws.Range("J9").Value = 123456781
$ws.Range("J10").Value = 123456782
$ws.Range("J11").Value = 123456783
$ws.Range("J12").Value = 123456784
$ws.Range("J13").Value = 123456785

# Reflect the final selection / active sheet state left behind by the edit.
$ws.Activate()
$ws.Range("J8").Select()

$wb.Save()
